# Update Name of Algo
# Apply updated values to the RandomForest result data sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = 13.5921
$ws.Range("B12").Value = 5.765099999999999
$ws.Range("E12").Value = 12.60079999999999
$ws.Range("E14").Value = 13.82910000000001
$ws.Range("E22").Value = 11.7041
